$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.2
$ws.Range("I2").Value = 3
$ws.Range("L2").Value = 1.34
$ws.Range("N2").Value = 3.45
$ws.Range("P2").Value = 1.86
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 1.72
$ws.Range("V2").Value = 1.51
$ws.Range("W2").Value = 1.46
$ws.Range("AJ2").Value = 55

# Row 3
$ws.Range("F3").Value = 2.64
$ws.Range("G3").Value = 3.3
$ws.Range("H3").Value = 2.6
$ws.Range("I3").Value = 3.2
$ws.Range("K3").Value = 3.85
$ws.Range("L3").Value = 1.38
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 3.35
$ws.Range("O3").Value = 1.27
$ws.Range("P3").Value = 1.91
$ws.Range("R3").Value = 1.32
$ws.Range("S3").Value = 2.86
$ws.Range("T3").Value = 1.05
$ws.Range("U3").Value = 1.05
$ws.Range("V3").Value = 1.45
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 14.5
$ws.Range("AB3").Value = 15
$ws.Range("AC3").Value = 9.800000000000001
$ws.Range("AD3").Value = 15
$ws.Range("AG3").Value = 15.5

# Row 4
$ws.Range("F4").Value = 2.36
$ws.Range("G4").Value = 2.92
$ws.Range("H4").Value = 2.84
$ws.Range("I4").Value = 3.6
$ws.Range("J4").Value = 3.15

# Row 5
$ws.Range("H5").Value = 1.42
$ws.Range("N5").Value = 5
$ws.Range("R5").Value = 1.54
$ws.Range("S5").Value = 2.74
$ws.Range("T5").Value = 1.95
$ws.Range("Z5").Value = 8.800000000000001
$ws.Range("AA5").Value = 12.5
$ws.Range("AD5").Value = 9.800000000000001
$ws.Range("AE5").Value = 15
$ws.Range("AJ5").Value = 300

# Row 6
$ws.Range("I6").Value = 1.91
$ws.Range("P6").Value = 2.3
$ws.Range("Q6").Value = 1.72
$ws.Range("S6").Value = 2.8
$ws.Range("Y6").Value = 11

# Row 8
$ws.Range("F8").Value = 2.34
$ws.Range("J8").Value = 3.45
$ws.Range("M8").Value = 1.08
$ws.Range("O8").Value = 1.34
$ws.Range("P8").Value = 1.9
$ws.Range("Q8").Value = 2.06
$ws.Range("R8").Value = 1.34
$ws.Range("S8").Value = 3.7
$ws.Range("T8").Value = 1.82
$ws.Range("U8").Value = 2.14
$ws.Range("X8").Value = 14.5
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 23
$ws.Range("AD8").Value = 15.5
$ws.Range("AF8").Value = 14.5
$ws.Range("AG8").Value = 11.5
$ws.Range("AH8").Value = 18.5
$ws.Range("AK8").Value = 26
$ws.Range("AL8").Value = 42
$ws.Range("AM8").Value = 120
$ws.Range("AN8").Value = 21
$ws.Range("AO8").Value = 44

# Row 9
$ws.Range("F9").Value = 12.5
$ws.Range("G9").Value = 13
$ws.Range("S9").Value = 2.88
$ws.Range("T9").Value = 2.3
$ws.Range("U9").Value = 1.73
$ws.Range("Y9").Value = 8.199999999999999
$ws.Range("Z9").Value = 7
$ws.Range("AG9").Value = 50
$ws.Range("AM9").Value = 440

# Row 10
$ws.Range("F10").Value = 2.72
$ws.Range("G10").Value = 2.74
$ws.Range("H10").Value = 2.88
$ws.Range("I10").Value = 2.9
$ws.Range("N10").Value = 3.95
$ws.Range("U10").Value = 2.26
$ws.Range("AA10").Value = 60
$ws.Range("AC10").Value = 7.6
$ws.Range("AE10").Value = 32
$ws.Range("AI10").Value = 55

# Row 11
$ws.Range("P11").Value = 2.36
$ws.Range("Q11").Value = 1.7
$ws.Range("S11").Value = 2.72
$ws.Range("U11").Value = 2.42
$ws.Range("AB11").Value = 24
$ws.Range("AE11").Value = 20
$ws.Range("AG11").Value = 17
$ws.Range("AO11").Value = 11.5

# Row 12
$ws.Range("F12").Value = 5
$ws.Range("H12").Value = 1.72
$ws.Range("I12").Value = 1.73
$ws.Range("N12").Value = 6.2
$ws.Range("P12").Value = 2.72
$ws.Range("S12").Value = 2.34
$ws.Range("T12").Value = 1.6
$ws.Range("U12").Value = 2.58
$ws.Range("Z12").Value = 13
$ws.Range("AC12").Value = 10.5
$ws.Range("AE12").Value = 15.5
$ws.Range("AF12").Value = 46
$ws.Range("AG12").Value = 20
